$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was updated
# from 45190 (2023-09-21) to 45192 (2023-09-23) for every data row
# (rows 2 through 110).
$ws.Range("C2:C110").Value = 45192
